$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "305.45"
$ws.Range("E2").Value = "-0.59%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "35.77"
$ws.Range("E3").Value = "-0.31%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.033"
$ws.Range("E4").Value = "-1.40%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07970"
$ws.Range("E5").Value = "-1.43%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "1.916"
$ws.Range("E6").Value = "-1.39%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.134"
$ws.Range("E7").Value = "-1.49%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.772"
$ws.Range("E8").Value = "0.32%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9194"
$ws.Range("E9").Value = "-0.91%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1299"
$ws.Range("E10").Value = "-5.70%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1910"
$ws.Range("E11").Value = "0.35%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09112"
$ws.Range("E12").Value = "-1.12%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03437"
$ws.Range("E13").Value = "0.84%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09835"
$ws.Range("E14").Value = "0.04%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001412"
$ws.Range("E15").Value = "-2.08%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006188"
$ws.Range("E16").Value = "6.80%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.730"
$ws.Range("E17").Value = "3.14%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.94%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.05%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").Value = "-2.32%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "5.161"
$ws.Range("E21").Value = "5.35%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2598"
$ws.Range("E22").Value = "6.20%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04418"
$ws.Range("E23").Value = "-0.33%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").Value = "0.71%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004632"
$ws.Range("E25").Value = "-4.03%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("E26").Value = "0.60%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004440"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01947"
$ws.Range("E39").Value = "-4.04%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05302"
$ws.Range("E40").Value = "7.77%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007598"
$ws.Range("E41").Value = "-0.17%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01018"
$ws.Range("E42").Value = "0.75%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1355"
$ws.Range("E43").Value = "-1.58%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002150"
$ws.Range("E44").Value = "2.17%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009955"
$ws.Range("E45").Value = "-9.63%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006074"
$ws.Range("E46").Value = "-5.88%"

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.26%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "64.96"
$ws.Range("E48").Value = "2.19%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "38.92%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.26%"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.26%"
